# Add support for a time.Time data type on the row struct.
#
# Adds two new columns to the "Hoja1" sheet:
#   H: "Date"   - a date/time value (numFmtId 22, "m/d/yyyy h:mm")
#   I: "Date 2" - the same date formatted/stored as plain text (numFmtId 49, "@")
# and widens the autofilter / used range / defined name to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header row ----------------------------------------------------------
$ws.Range("H1").Value = "Date"
$ws.Range("I1").Value = "Date 2"

# --- Data rows 2-7: date/time value + its text representation ------------
$dateValue = Get-Date -Year 2022 -Month 9 -Day 20 -Hour 13 -Minute 33 -Second 1

for ($r = 2; $r -le 7; $r++) {
    $hCell = $ws.Range("H$r")
    $hCell.Font.Underline = $false          # keep default font (fontId 0), even on rows with custom row format
    $hCell.NumberFormat = "m/d/yy h:mm"     # normalizes to builtin numFmtId 22 (m/d/yyyy h:mm)
    $hCell.Value = $dateValue

    $iCell = $ws.Range("I$r")
    $iCell.Font.Underline = $false          # keep default font (fontId 0)
    $iCell.NumberFormat = "@"               # builtin numFmtId 49, forces text storage
    $iCell.Value = "2022-09-20"
}

# --- Column widths for the two new columns --------------------------------
$ws.Columns.Item(8).ColumnWidth = 19.3
$ws.Columns.Item(9).ColumnWidth = 19.5

# --- Extend the autofilter to cover the full new range A1:H7 -------------
$ws.Range("A1:H7").AutoFilter() | Out-Null

# --- Reset the selection back to the top-left cell ------------------------
$ws.Range("A1").Select()
